# Auto-generated edit script applying the crypto price/volume update diff
# (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.534.41'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.098.66'
$ws.Range('E3').Value = '  -0.69%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.65'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.11'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.091.23'
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.160'
$ws.Range('E10').Value = '  +6.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.61'
$ws.Range('E11').Value = '  -3.09%  '
$ws.Range('E12').Value = '  -2.35%  '
$ws.Range('E13').Value = '  -1.01%  '
$ws.Range('E14').Value = '  +4.17%  '
$ws.Range('E15').Value = '  -1.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.612.41'
$ws.Range('E16').Value = '  -0.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.369.28'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.06'
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.092.25'
$ws.Range('E19').Value = '  -0.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '460.45'
$ws.Range('E20').Value = '  -1.54%  '
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('E22').Value = '  -0.48%  '
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.89'
$ws.Range('E24').Value = '  -3.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.04'
$ws.Range('E25').Value = '  -1.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.21'
$ws.Range('E26').Value = '  +1.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.17'
$ws.Range('E28').Value = '  +9.85%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  -0.59%  '
$ws.Range('E31').Value = '  -1.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.91'
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.109'
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.59'
$ws.Range('E34').Value = '  -1.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0846'
$ws.Range('E35').Value = '  -2.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.37'
$ws.Range('E36').Value = '  +1.95%  '
$ws.Range('E37').Value = '  -0.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.30'
$ws.Range('E38').Value = '  -5.06%  '
$ws.Range('E39').Value = '  -0.79%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '50.21'
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '435.57'
$ws.Range('E41').Value = '  +0.76%  '
$ws.Range('E42').Value = '  -0.68%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.881.14'
$ws.Range('E43').Value = '  -1.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0367'
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.272'
$ws.Range('E45').Value = '  -2.28%  '
$ws.Range('E46').Value = '  -3.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '36.35'
$ws.Range('E47').Value = '  +3.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.69'
$ws.Range('E48').Value = '  +0.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.110'
$ws.Range('E50').Value = '  -1.38%  '
$ws.Range('E51').Value = '  -2.37%  '
